# Auto-generated Excel COM-interop script
# Applies crypto price/volume updates from the commit diff, including
# a row swap between InjectiveProtocol (row 47->48) and Mantle (row 48->47).
#
# Note: several Price (column D) values look like plain decimals (e.g. '543.53').
# The source workbook stores Price/Volume as literal TEXT cells (t="inlineStr"),
# not numbers. A bare Range.Value assignment of such a string would be
# auto-coerced to a numeric cell by Excel, which would NOT match the original
# text-cell representation. To force text while keeping the cell's original
# (default/general) style, we assign the value with a leading apostrophe
# (Excel's standard 'treat as text' marker, stored, not displayed) and then
# reset Style to 'Normal' to drop the transient quote-prefix style Excel
# attaches to the cell -- leaving the cell looking exactly as it did before,
# just holding the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.518.03'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '2.331.27'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'543.53"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.74%  '
$ws.Range('D6').Value = "'135.20"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('D7').Value = "'0.992"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('D9').Value = '2.365.99'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('E10').Value = '  +1.43%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('D12').Value = "'5.40"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('D13').Value = "'0.354"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.16%  '
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').Value = '2.755.48'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '57.641.19'
$ws.Range('E16').Value = '  +1.67%  '
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '2.339.62'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = "'338.28"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.68%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('E21').Value = '  +0.80%  '
$ws.Range('D22').Value = "'6.85"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.69%  '
$ws.Range('D23').Value = "'0.998"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = "'61.80"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.05%  '
$ws.Range('E25').Value = '  +2.57%  '
$ws.Range('E26').Value = '  -2.66%  '
$ws.Range('D27').Value = "'0.997"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +5.30%  '
$ws.Range('D29').Value = "'174.75"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.83%  '
$ws.Range('D30').Value = "'1.76"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.95%  '
$ws.Range('D31').Value = '0.0₃0738'
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('D32').Value = "'6.16"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('D33').Value = "'18.56"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('E35').Value = '  +12.64%  '
$ws.Range('D36').Value = "'0.990"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').Value = "'1.26"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('D38').Value = "'4.13"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.69%  '
$ws.Range('E39').Value = '  +2.70%  '
$ws.Range('D40').Value = "'39.46"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('D41').Value = "'149.23"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('D43').Value = "'3.64"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.40%  '
$ws.Range('D44').Value = "'284.94"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = "'0.563"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = "'18.76"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.55%  '
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').Value = "'17.56"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('E51').Value = '  +8.31%  '
